# Apply crypto price/volume updates to the active sheet, row-by-row,
# matching the committed diff exactly (including the two 13/14 and 34/35 coin swaps).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.725.70"
$ws.Range("E2").Value = "  +0.18%  "
$ws.Range("D3").Value = "1.647.88"
$ws.Range("E3").Value = "  +0.61%  "
$ws.Range("E4").Value = "  +0.30%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "215.99"
$ws.Range("E5").Value = "  +1.31%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.503"
$ws.Range("E6").Value = "  +0.45%  "
$ws.Range("E7").Value = "  +0.31%  "
$ws.Range("E8").Value = "  -1.05%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.0627"
$ws.Range("E9").Value = "  +0.53%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "19.31"
$ws.Range("E10").Value = "  +0.46%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0843"
$ws.Range("E11").Value = "  +0.13%  "
$ws.Range("D12").Value = "1.879.06"
$ws.Range("E12").Value = "  +0.69%  "
$ws.Range("B13").Value = "Polkadot"
$ws.Range("C13").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "4.22"
$ws.Range("E13").Value = "  +2.76%  "
$ws.Range("B14").Value = "WrappedEther"
$ws.Range("C14").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D14").Value = "1.645.16"
$ws.Range("E14").Value = "  +0.69%  "
$ws.Range("E15").Value = "  +0.97%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "66.39"
$ws.Range("E16").Value = "  +4.86%  "
$ws.Range("D17").Value = "26.765.05"
$ws.Range("E17").Value = "  +0.29%  "
$ws.Range("D18").Value = "0.0₃0752"
$ws.Range("E18").Value = "  +1.07%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "220.04"
$ws.Range("E19").Value = "  +0.73%  "
$ws.Range("E20").Value = "  +0.24%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "4.39"
$ws.Range("E21").Value = "  +1.55%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.33"
$ws.Range("E22").Value = "  +1.89%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "9.54"
$ws.Range("E23").Value = "  +0.62%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.12"
$ws.Range("E24").Value = "  +10.54%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "146.79"
$ws.Range("E25").Value = "  -1.22%  "
$ws.Range("E26").Value = "  +0.21%  "
$ws.Range("E27").Value = "  +0.14%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "7.09"
$ws.Range("E28").Value = "  +3.22%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "15.87"
$ws.Range("E29").Value = "  +2.70%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.0518"
$ws.Range("E30").Value = "  +1.45%  "
$ws.Range("E31").Value = "  +0.51%  "
$ws.Range("E32").Value = "  +2.92%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.04"
$ws.Range("E33").Value = "  +3.09%  "
$ws.Range("B34").Value = "Maker"
$ws.Range("C34").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D34").Value = "1.286.76"
$ws.Range("E34").Value = "  +7.26%  "
$ws.Range("B35").Value = "LidoDAOToken"
$ws.Range("C35").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.56"
$ws.Range("E35").Value = "  +2.76%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.0184"
$ws.Range("E36").Value = "  +6.14%  "
$ws.Range("E37").Value = "  +0.81%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.828"
$ws.Range("E38").Value = "  +1.98%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.523"
$ws.Range("E39").Value = "  +3.22%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.01"
$ws.Range("E40").Value = "  +0.33%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.811"
$ws.Range("E41").Value = "  +2.17%  "
$ws.Range("E42").Value = "  -2.77%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "5.43"
$ws.Range("E43").Value = "  +0.24%  "
$ws.Range("D44").Value = "1.789.35"
$ws.Range("E44").Value = "  +0.86%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "93.54"
$ws.Range("E45").Value = "  +1.33%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "59.35"
$ws.Range("E46").Value = "  +8.08%  "
$ws.Range("E47").Value = "  +3.72%  "
$ws.Range("E48").Value = "  +0.92%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "7.76"
$ws.Range("E49").Value = "  +1.66%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0974"
$ws.Range("E50").Value = "  +2.75%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.408"
$ws.Range("E51").Value = "  -0.59%  "
